{"js": "// Replace the two-digit multiplication problems in the table with the\n// newly generated set of problems. Each old equation string is unique\n// within the document, so a simple exact search-and-replace per pair\n// is sufficient and keeps the original run formatting intact.\nconst replacements = [\n  [\"18\u00d768=\", \"49\u00d777=\"],\n  [\"68\u00d798=\", \"96\u00d754=\"],\n  [\"49\u00d769=\", \"89\u00d717=\"],\n  [\"51\u00d799=\", \"77\u00d768=\"],\n  [\"34\u00d733=\", \"91\u00d748=\"],\n  [\"17\u00d781=\", \"77\u00d721=\"],\n  [\"65\u00d739=\", \"46\u00d761=\"],\n  [\"68\u00d731=\", \"37\u00d781=\"],\n  [\"19\u00d735=\", \"94\u00d737=\"],\n  [\"18\u00d760=\", \"30\u00d720=\"],\n  [\"18\u00d725=\", \"32\u00d719=\"],\n  [\"26\u00d727=\", \"46\u00d788=\"],\n  [\"55\u00d741=\", \"32\u00d744=\"],\n  [\"35\u00d758=\", \"88\u00d774=\"],\n  [\"52\u00d723=\", \"77\u00d783=\"],\n  [\"42\u00d789=\", \"55\u00d761=\"],\n  [\"40\u00d752=\", \"43\u00d714=\"],\n  [\"89\u00d781=\", \"90\u00d774=\"],\n  [\"65\u00d734=\", \"84\u00d768=\"],\n  [\"68\u00d718=\", \"36\u00d763=\"],\n  [\"48\u00d732=\", \"19\u00d792=\"],\n  [\"38\u00d765=\", \"65\u00d785=\"],\n  [\"18\u00d766=\", \"89\u00d733=\"],\n  [\"18\u00d747=\", \"20\u00d776=\"],\n  [\"43\u00d715=\", \"45\u00d794=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the table with the\n# newly generated set of problems. Each old equation string is unique\n# within the document, so a Find/Replace pass per pair (ReplaceAll)\n# is sufficient and preserves the original run formatting.\n\n$d = $word.ActiveDocument\n\n$pairs = [ordered]@{\n    \"18\u00d768=\" = \"49\u00d777=\"\n    \"68\u00d798=\" = \"96\u00d754=\"\n    \"49\u00d769=\" = \"89\u00d717=\"\n    \"51\u00d799=\" = \"77\u00d768=\"\n    \"34\u00d733=\" = \"91\u00d748=\"\n    \"17\u00d781=\" = \"77\u00d721=\"\n    \"65\u00d739=\" = \"46\u00d761=\"\n    \"68\u00d731=\" = \"37\u00d781=\"\n    \"19\u00d735=\" = \"94\u00d737=\"\n    \"18\u00d760=\" = \"30\u00d720=\"\n    \"18\u00d725=\" = \"32\u00d719=\"\n    \"26\u00d727=\" = \"46\u00d788=\"\n    \"55\u00d741=\" = \"32\u00d744=\"\n    \"35\u00d758=\" = \"88\u00d774=\"\n    \"52\u00d723=\" = \"77\u00d783=\"\n    \"42\u00d789=\" = \"55\u00d761=\"\n    \"40\u00d752=\" = \"43\u00d714=\"\n    \"89\u00d781=\" = \"90\u00d774=\"\n    \"65\u00d734=\" = \"84\u00d768=\"\n    \"68\u00d718=\" = \"36\u00d763=\"\n    \"48\u00d732=\" = \"19\u00d792=\"\n    \"38\u00d765=\" = \"65\u00d785=\"\n    \"18\u00d766=\" = \"89\u00d733=\"\n    \"18\u00d747=\" = \"20\u00d776=\"\n    \"43\u00d715=\" = \"45\u00d794=\"\n}\n\nforeach ($oldText in $pairs.Keys) {\n    $newText = $pairs[$oldText]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
